$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "transactions"
$ws.Range("A2").Value = "jsaf123jlsa513j"
$ws.Range("B2").Value = "SLA breached for HITBTC tech issue"
$ws.Range("C2").Value = 11
$ws.Range("A3").Value = "vbzsdfgdsg1234"
$ws.Range("B3").Value = "SLA breached for binance tech issue"
$ws.Range("C3").Value = 12
$ws.Range("B1").Value = "comment"

$ws.Columns.Item(1).ColumnWidth = 26.022135416666668
$ws.Columns.Item(2).ColumnWidth = 31.592447916666668

$ws.Range("E4").Select()
